# Fruta / hortaliza, semanal
# Insert one new weekly data row before row 332 (shifting the existing
# rows 332-360 down to 333-361) and populate the new row with the latest
# week's reading for Femacal de La Calera - Arándano (blue).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 332:360 down by one row, creating a blank row 332.
$ws.Rows.Item(332).Insert()

# Fill the newly inserted row 332 with this week's data.
$ws.Cells.Item(332, 1).Value = 3
$ws.Cells.Item(332, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(332, 3).Value = "Coquimbo"
$ws.Cells.Item(332, 4).Value = 45154
$ws.Cells.Item(332, 5).Value = 5
$ws.Cells.Item(332, 6).Value = "Fruta"
$ws.Cells.Item(332, 7).Value = 100101
$ws.Cells.Item(332, 8).Value = "Berries"
$ws.Cells.Item(332, 9).Value = 100101001
$ws.Cells.Item(332, 10).Value = "Arándano (blue)"
$ws.Cells.Item(332, 11).Value = "Sin especificar"
$ws.Cells.Item(332, 12).Value = "Primera"
$ws.Cells.Item(332, 13).Value = 30
$ws.Cells.Item(332, 14).Value = 13000
$ws.Cells.Item(332, 15).Value = 13000
$ws.Cells.Item(332, 16).Value = 13000
$ws.Cells.Item(332, 17).Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(332, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(332, 19).Value = 8667
$ws.Cells.Item(332, 20).Value = 1.5
